$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the quote/document date (A1) from 24-Apr-2024 (45406) to 24-May-2024 (45436)
$ws.Range("A1").Value = 45436

# Update price list values
$ws.Range("D35").Value = 32935
$ws.Range("D36").Value = 7128
